$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wrap the hotel-name link template (C3) in a bootstrap column div
$ws.Range("C3").Value = '<div class="col-md-4 col-sm-4 no-padding"><a class="hotel-name" href="#LINK">#TEXT</a></div>'
